$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Barcode column (K) becomes numeric quantity-like values instead of text barcodes,
# and a new "Quantity" column (O) mirrors the same numeric values.
$ws.Range("K2").Value = 1000
$ws.Range("O2").Value = 1000

$ws.Range("K3").Value = 1001
$ws.Range("O3").Value = 1001

$ws.Range("K4").Value = 1002
$ws.Range("O4").Value = 1002

$ws.Range("K5").Value = 1003
$ws.Range("O5").Value = 1003

# Add the new "Quantity" header in column O, matching the header style used
# by the neighbouring header cells (N1). PasteSpecial(xlPasteFormats) only
# copies the formatting (bold/italic header style), leaving the value we
# already set above untouched.
$ws.Range("O1").Value = "Quantity"
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)

# Widen column K (Barcode) to fit the new numeric values.
$ws.Columns.Item(11).ColumnWidth = 42.5

# Update the active selection to O6 to mirror the new last column.
$ws.Range("O6").Select()
